# Atualiza a coluna (P)rioridade e (B)aseline de características específicas
# na tabela "Lista de Características (Prioridade x Esforço x Risco x Baseline)".
#
# A tabela tem 6 colunas: #, Características, (P), (E), (R), (B)
# A linha 1 da tabela é o cabeçalho, então a característica de número N
# corresponde à linha (N + 1) da tabela do Word.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellValue($table, $row, $col, $newValue) {
    $cell = $table.Cell($row, $col)
    $rng = $cell.Range
    # Exclui a marca de fim de célula / quebra de parágrafo do intervalo.
    $rng.End = $rng.End - 1
    $rng.Text = $newValue
}

# Coluna (P): de "C" para "I"
Set-CellValue $t 5  3 "I"
Set-CellValue $t 7  3 "I"
Set-CellValue $t 8  3 "I"
Set-CellValue $t 9  3 "I"
Set-CellValue $t 11 3 "I"
Set-CellValue $t 12 3 "I"
Set-CellValue $t 14 3 "I"
Set-CellValue $t 20 3 "I"

# Coluna (P): de "C" para "U", e coluna (B): de "1" para "2"
Set-CellValue $t 18 3 "U"
Set-CellValue $t 18 6 "2"
Set-CellValue $t 19 3 "U"
Set-CellValue $t 19 6 "2"
Set-CellValue $t 21 3 "U"
Set-CellValue $t 21 6 "2"
Set-CellValue $t 22 3 "U"
Set-CellValue $t 22 6 "2"
